# Apply cryptocurrency price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.280.65"
$ws.Range("E2").Value = "  +6.07%  "

$ws.Range("D3").Value = "3.120.30"
$ws.Range("E3").Value = "  +3.83%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'585.44"
$ws.Range("E5").Value = "  +3.86%  "

$ws.Range("D6").Value = "'145.40"
$ws.Range("E6").Value = "  +4.15%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.111.58"
$ws.Range("E8").Value = "  +3.94%  "

$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  +12.78%  "

$ws.Range("D11").Value = "'5.80"
$ws.Range("E11").Value = "  +10.01%  "

$ws.Range("E12").Value = "  +2.97%  "

$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +7.95%  "

$ws.Range("D14").Value = "'35.63"
$ws.Range("E14").Value = "  +4.91%  "

$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "3.635.75"
$ws.Range("E16").Value = "  +3.76%  "

$ws.Range("D18").Value = "63.189.08"
$ws.Range("E18").Value = "  +5.93%  "

$ws.Range("D19").Value = "3.115.62"
$ws.Range("E19").Value = "  +3.64%  "

$ws.Range("D20").Value = "'466.12"
$ws.Range("E20").Value = "  +6.11%  "

$ws.Range("D21").Value = "'14.06"
$ws.Range("E21").Value = "  +3.12%  "

$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").Value = "'7.53"
$ws.Range("E23").Value = "  +5.97%  "

$ws.Range("E24").Value = "  -1.69%  "

$ws.Range("D25").Value = "'82.08"
$ws.Range("E25").Value = "  +2.00%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  +0.68%  "

$ws.Range("E28").Value = "  +5.00%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'8.28"
$ws.Range("E29").Value = "  +6.46%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("E31").Value = "  +8.73%  "

$ws.Range("D32").Value = "'27.00"
$ws.Range("E32").Value = "  +4.29%  "

$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = "  +3.84%  "

$ws.Range("D34").Value = "0.0₃0872"
$ws.Range("E34").Value = "  +11.87%  "

$ws.Range("D35").Value = "'2.42"
$ws.Range("E35").Value = "  +15.41%  "

$ws.Range("E36").Value = "  +5.11%  "

$ws.Range("E37").Value = "  +20.29%  "

$ws.Range("D38").Value = "'6.05"
$ws.Range("E38").Value = "  +2.13%  "

$ws.Range("D39").Value = "'50.86"
$ws.Range("E39").Value = "  +3.61%  "

$ws.Range("D40").Value = "'432.60"
$ws.Range("E40").Value = "  +7.52%  "

$ws.Range("D41").Value = "'8.72"
$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("D42").Value = "2.911.43"
$ws.Range("E42").Value = "  +5.16%  "

$ws.Range("E43").Value = "  +4.12%  "

$ws.Range("E44").Value = "  +10.38%  "

$ws.Range("E45").Value = "  +5.27%  "

$ws.Range("D46").Value = "'2.17"
$ws.Range("E46").Value = "  +6.29%  "

$ws.Range("D47").Value = "'35.18"
$ws.Range("E47").Value = "  +2.60%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").Value = "'123.77"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("E50").Value = "  +0.59%  "

$ws.Range("D51").Value = "'24.54"
$ws.Range("E51").Value = "  +3.50%  "
